$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
for ($i=0; $i -lt $wb.Styles.Count; $i++) {
    Write-Output $wb.Styles.Item($i+1).Name
}
